# Rename the "*img" sheets to "img*" and make the last one (now "imge")
# the active tab, matching the source workbook's state.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

$wb.Worksheets.Item("imge").Activate()
